$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete column H ("Anthony's desktop") entirely - no longer needed
$ws.Columns.Item(8).Delete()

# Row 6 (BMsolveGSSA): change source from Daryl's laptop (F) to Kerk's home (D)
$ws.Range("F6").ClearContents()
$ws.Range("D6").Value = 29.671941295771202
$ws.Range("B6").Formula = "=D6/D2"

# Row 7 (BMsolveVFI): add new result from Kerk's laptop (C)
$ws.Range("C7").Value = 93.098098153946907
$ws.Range("B7").Formula = "=C7/C2"

# Row 11 (BMsimVFI): add new result from Kerk's laptop (C)
$ws.Range("C11").Value = 4624.24204680603
$ws.Range("B11").Formula = "=C11/C2"

# Update view state: move active selection to B11 (also repositions the frozen pane)
$ws.Range("B11").Select()
